$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.481.76"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.909.53"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.66%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "325.70"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4852"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4063"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08171"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.013"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "23.43"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.74%  "
$ws.Range("D12").Value = "1.918.95"
$ws.Range("E12").Value = "  -0.97%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.017"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.166"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "90.42"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.008"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.06777"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.42%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001037"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "29.499.88"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.631"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.10%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.76"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "2.167.66"
$ws.Range("E25").Value = "  +0.19%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.12"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.84%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "6.522"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +8.42%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "20.09"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.114"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "120.43"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.29%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.030"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.77%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09529"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.514"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.89%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.561"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.52%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.391"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02275"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.36%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06119"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.179"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "10.86"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.62%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5961"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.23%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "7.974"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.58%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1855"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -4.90%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "12.50"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.84%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.07617"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.25%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5577"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.951"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +3.08%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "72.55"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.89%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.406"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.49%  "
